$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell 'D2' '42.773.43'
Set-TextCell 'E2' '  +0.51%  '

Set-TextCell 'D3' '2.523.41'
Set-TextCell 'E3' '  +0.26%  '

Set-TextCell 'E4' '  +0.04%  '

Set-TextCell 'D5' '313.96'
Set-TextCell 'E5' '  +1.71%  '

Set-TextCell 'D6' '95.94'
Set-TextCell 'E6' '  -0.67%  '

Set-TextCell 'D7' '0.576'
Set-TextCell 'E7' '  -1.69%  '

Set-TextCell 'E8' '  -0.08%  '

Set-TextCell 'D9' '0.534'
Set-TextCell 'E9' '  -1.07%  '

Set-TextCell 'D10' '36.13'
Set-TextCell 'E10' '  -1.67%  '

Set-TextCell 'D11' '0.0810'
Set-TextCell 'E11' '  -0.47%  '

Set-TextCell 'D12' '7.55'
Set-TextCell 'E12' '  -2.55%  '

Set-TextCell 'D13' '0.109'
Set-TextCell 'E13' '  -3.73%  '

Set-TextCell 'D14' '2.910.11'
Set-TextCell 'E14' '  +0.32%  '

Set-TextCell 'D15' '2.545.64'
Set-TextCell 'E15' '  +1.56%  '

Set-TextCell 'D16' '15.26'
Set-TextCell 'E16' '  -3.58%  '

Set-TextCell 'D17' '0.854'
Set-TextCell 'E17' '  -1.12%  '

Set-TextCell 'D18' '42.827.51'
Set-TextCell 'E18' '  +0.72%  '

Set-TextCell 'B19' 'InternetComputer(DFINITY)'
Set-TextCell 'C19' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D19' '12.89'
Set-TextCell 'E19' '  -0.82%  '

Set-TextCell 'B20' 'Uniswap'
Set-TextCell 'C20' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D20' '6.74'
Set-TextCell 'E20' '  +4.20%  '

Set-TextCell 'D21' '0.0₃0962'
Set-TextCell 'E21' '  -1.29%  '

Set-TextCell 'D22' '69.68'
Set-TextCell 'E22' '  -2.60%  '

Set-TextCell 'D23' '254.56'
Set-TextCell 'E23' '  +0.23%  '

Set-TextCell 'E24' '  +0.13%  '

Set-TextCell 'E25' '  +1.83%  '

Set-TextCell 'D26' '26.63'
Set-TextCell 'E26' '  -1.65%  '

Set-TextCell 'E27' '  -0.12%  '

Set-TextCell 'E28' '  +3.63%  '

Set-TextCell 'D29' '41.03'
Set-TextCell 'E29' '  +9.07%  '

Set-TextCell 'D30' '10.37'
Set-TextCell 'E30' '  +1.80%  '

Set-TextCell 'D31' '5.95'
Set-TextCell 'E31' '  -0.13%  '

Set-TextCell 'D32' '157.80'
Set-TextCell 'E32' '  +2.53%  '

Set-TextCell 'D33' '19.56'
Set-TextCell 'E33' '  +2.05%  '

Set-TextCell 'E34' '  +3.53%  '

Set-TextCell 'E35' '  +2.86%  '

Set-TextCell 'E36' '  +0.66%  '

Set-TextCell 'D37' '0.0781'
Set-TextCell 'E37' '  -1.00%  '

Set-TextCell 'D38' '0.112'
Set-TextCell 'E38' '  -2.37%  '

Set-TextCell 'D39' '0.119'
Set-TextCell 'E39' '  -0.96%  '

Set-TextCell 'D40' '23.52'
Set-TextCell 'E40' '  -3.61%  '

Set-TextCell 'D41' '2.32'
Set-TextCell 'E41' '  +13.60%  '

Set-TextCell 'D42' '0.0305'
Set-TextCell 'E42' '  +0.86%  '

Set-TextCell 'B43' 'NEARProtocol'
Set-TextCell 'C43' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D43' '3.33'
Set-TextCell 'E43' '  -2.18%  '

Set-TextCell 'D44' '3.80'
Set-TextCell 'E44' '  -2.12%  '

Set-TextCell 'B45' 'FirstDigitalUSD'
Set-TextCell 'C45' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 'D45' '1.00'
Set-TextCell 'E45' '  +0.40%  '

Set-TextCell 'D46' '2.044.44'
Set-TextCell 'E46' '  +0.22%  '

Set-TextCell 'D47' '85.15'
Set-TextCell 'E47' '  +0.44%  '

Set-TextCell 'D48' '8.92'
Set-TextCell 'E48' '  -0.83%  '

Set-TextCell 'D49' '75.82'
Set-TextCell 'E49' '  +3.85%  '

Set-TextCell 'D50' '106.27'
Set-TextCell 'E50' '  +4.54%  '

Set-TextCell 'D51' '2.762.50'
Set-TextCell 'E51' '  +0.18%  '
